$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (43 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2983.8667
$ws.Range("I15").Value = 2983.8667
$ws.Range("K15").Value = 8951.6001
$ws.Range("M15").Value = -8782.6001
$ws.Range("H76").Value = 2940.5386
$ws.Range("I76").Value = 2846.5557
$ws.Range("J76").Value = 3152
$ws.Range("K76").Value = 2846.5557
$ws.Range("L76").Value = 3152
$ws.Range("M76").Value = -2531.5557
$ws.Range("N76").Value = -3782
$ws.Range("H79").Value = 2940.5386
$ws.Range("I79").Value = 2846.5557
$ws.Range("J79").Value = 3152
$ws.Range("K79").Value = 2846.5557
$ws.Range("L79").Value = 3152
$ws.Range("M79").Value = -1754.5557
$ws.Range("N79").Value = -5336
$ws.Range("H132").Value = 2987383.8
$ws.Range("I132").Value = 3510900.5
$ws.Range("J132").Value = 3338.1
$ws.Range("K132").Value = 10532701.5
$ws.Range("L132").Value = 10014.3
$ws.Range("N132").Value = -15074.3
$ws.Range("M132").Value = -10530171.5
$ws.Range("H137").Value = 3066.6606
$ws.Range("I137").Value = 3416.0256
$ws.Range("K137").Value = 10248.0768
$ws.Range("M137").Value = -7698.076799999999
$ws.Range("H138").Value = 6126.174
$ws.Range("I138").Value = 1907.3158
$ws.Range("J138").Value = 9095
$ws.Range("K138").Value = 5721.9474
$ws.Range("L138").Value = 27285
$ws.Range("M138").Value = -581.9474
$ws.Range("N138").Value = -37565
$ws.Range("H141").Value = 496772.2
$ws.Range("I141").Value = 4005.725
$ws.Range("J141").Value = 2139327
$ws.Range("K141").Value = 12017.175
$ws.Range("L141").Value = 6417981
$ws.Range("M141").Value = -6837.174999999999
$ws.Range("N141").Value = -6428341

# --- Sheet: ARM (53 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4135.24
$ws.Range("I32").Value = 3040.9768
$ws.Range("J32").Value = 10857.143
$ws.Range("K32").Value = 3040.9768
$ws.Range("L32").Value = 10857.143
$ws.Range("M32").Value = -2753.9768
$ws.Range("N32").Value = -11431.143
$ws.Range("H61").Value = 1152.7174
$ws.Range("I61").Value = 643.4524
$ws.Range("J61").Value = 6500
$ws.Range("K61").Value = 643.4524
$ws.Range("L61").Value = 6500
$ws.Range("M61").Value = -431.4524
$ws.Range("N61").Value = -6924
$ws.Range("H74").Value = 974.5417
$ws.Range("I74").Value = 790.4091
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 790.4091
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = 83.59090000000003
$ws.Range("N74").Value = -4748
$ws.Range("H77").Value = 974.5417
$ws.Range("I77").Value = 790.4091
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 3952.0455
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = 415.9545000000003
$ws.Range("N77").Value = -23736
$ws.Range("H122").Value = 1537.6604
$ws.Range("I122").Value = 1139.0975
$ws.Range("J122").Value = 2899.4167
$ws.Range("K122").Value = 3417.2925
$ws.Range("L122").Value = 8698.250100000001
$ws.Range("M122").Value = -967.2925000000005
$ws.Range("N122").Value = -13598.2501
$ws.Range("H132").Value = 2422.7021
$ws.Range("I132").Value = 1725.1892
$ws.Range("J132").Value = 5003.5
$ws.Range("K132").Value = 5175.5676
$ws.Range("L132").Value = 15010.5
$ws.Range("M132").Value = -2645.5676
$ws.Range("N132").Value = -20070.5
$ws.Range("H136").Value = 1152.7174
$ws.Range("I136").Value = 643.4524
$ws.Range("J136").Value = 6500
$ws.Range("K136").Value = 1930.3572
$ws.Range("L136").Value = 19500
$ws.Range("M136").Value = 619.6428000000001
$ws.Range("N136").Value = -24600
$ws.Range("H138").Value = 74250
$ws.Range("J138").Value = 74250
$ws.Range("L138").Value = 74250
$ws.Range("N138").Value = -84530

# --- Sheet: BSM (25 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 684729.0600000001
$ws.Range("I86").Value = 1015550.3
$ws.Range("J86").Value = 59844.445
$ws.Range("K86").Value = 1015550.3
$ws.Range("L86").Value = 59844.445
$ws.Range("M86").Value = -1014427.3
$ws.Range("N86").Value = -62090.445
$ws.Range("H89").Value = 684729.0600000001
$ws.Range("I89").Value = 1015550.3
$ws.Range("J89").Value = 59844.445
$ws.Range("K89").Value = 5077751.5
$ws.Range("L89").Value = 299222.225
$ws.Range("M89").Value = -5072135.5
$ws.Range("N89").Value = -310454.225
$ws.Range("H107").Value = 1923.1791
$ws.Range("I107").Value = 1579.6875
$ws.Range("J107").Value = 2790.9473
$ws.Range("K107").Value = 1579.6875
$ws.Range("L107").Value = 2790.9473
$ws.Range("M107").Value = 340.3125
$ws.Range("N107").Value = -6630.9473
$ws.Range("H134").Value = 3973.1785
$ws.Range("I134").Value = 3502.4285
$ws.Range("K134").Value = 10507.2855
$ws.Range("M134").Value = -7972.2855

# --- Sheet: CRP (43 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8067201
$ws.Range("I58").Value = 1663.7556
$ws.Range("J58").Value = 29417152
$ws.Range("K58").Value = 1663.7556
$ws.Range("L58").Value = 29417152
$ws.Range("M58").Value = -1460.7556
$ws.Range("N58").Value = -29417558
$ws.Range("H86").Value = 5649.5625
$ws.Range("I86").Value = 4343.4287
$ws.Range("J86").Value = 6665.4443
$ws.Range("K86").Value = 4343.4287
$ws.Range("L86").Value = 6665.4443
$ws.Range("M86").Value = -3220.4287
$ws.Range("N86").Value = -8911.444299999999
$ws.Range("H89").Value = 5649.5625
$ws.Range("I89").Value = 4343.4287
$ws.Range("J89").Value = 6665.4443
$ws.Range("K89").Value = 21717.1435
$ws.Range("L89").Value = 33327.2215
$ws.Range("M89").Value = -16101.1435
$ws.Range("N89").Value = -44559.2215
$ws.Range("H132").Value = 2105.422
$ws.Range("I132").Value = 1611.1765
$ws.Range("J132").Value = 3633.0908
$ws.Range("K132").Value = 4833.529500000001
$ws.Range("L132").Value = 10899.2724
$ws.Range("M132").Value = -2303.529500000001
$ws.Range("N132").Value = -15959.2724
$ws.Range("H134").Value = 1749.7028
$ws.Range("I134").Value = 953.35486
$ws.Range("K134").Value = 2860.06458
$ws.Range("M134").Value = -325.0645800000002
$ws.Range("H136").Value = 8067201
$ws.Range("I136").Value = 1663.7556
$ws.Range("J136").Value = 29417152
$ws.Range("K136").Value = 4991.266799999999
$ws.Range("L136").Value = 88251456
$ws.Range("M136").Value = -2441.266799999999
$ws.Range("N136").Value = -88256556
$ws.Range("H141").Value = 24888.889
$ws.Range("J141").Value = 24888.889
$ws.Range("L141").Value = 24888.889
$ws.Range("N141").Value = -35248.889

# --- Sheet: CUL (14 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 125.454544
$ws.Range("I33").Value = 76
$ws.Range("J33").Value = 166.66667
$ws.Range("K33").Value = 456
$ws.Range("L33").Value = 1000.00002
$ws.Range("M33").Value = -173
$ws.Range("N33").Value = -1566.00002
$ws.Range("H137").Value = 2943.7585
$ws.Range("I137").Value = 2691.9167
$ws.Range("J137").Value = 3121.5293
$ws.Range("K137").Value = 8075.750100000001
$ws.Range("L137").Value = 9364.5879
$ws.Range("M137").Value = -2975.750100000001
$ws.Range("N137").Value = -19564.5879

# --- Sheet: GSM (14 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3217.4736
$ws.Range("I122").Value = 2190.24
$ws.Range("J122").Value = 5192.923
$ws.Range("K122").Value = 6570.719999999999
$ws.Range("L122").Value = 15578.769
$ws.Range("M122").Value = -4120.719999999999
$ws.Range("N122").Value = -20478.769
$ws.Range("H132").Value = 3691.4473
$ws.Range("I132").Value = 3454.7036
$ws.Range("J132").Value = 4272.5454
$ws.Range("K132").Value = 10364.1108
$ws.Range("L132").Value = 12817.6362
$ws.Range("M132").Value = -7834.110799999999
$ws.Range("N132").Value = -17877.6362

# --- Sheet: LTW (18 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2775.111
$ws.Range("I122").Value = 2516.16
$ws.Range("K122").Value = 7548.48
$ws.Range("M122").Value = -5098.48
$ws.Range("H132").Value = 1830.9207
$ws.Range("I132").Value = 1162.8298
$ws.Range("J132").Value = 3793.4375
$ws.Range("K132").Value = 3488.4894
$ws.Range("L132").Value = 11380.3125
$ws.Range("M132").Value = -958.4893999999999
$ws.Range("N132").Value = -16440.3125
$ws.Range("H136").Value = 1660.0392
$ws.Range("I136").Value = 1151.4318
$ws.Range("J136").Value = 4857
$ws.Range("K136").Value = 3454.2954
$ws.Range("L136").Value = 14571
$ws.Range("M136").Value = -904.2954
$ws.Range("N136").Value = -19671

# --- Sheet: WVR (21 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2565534
$ws.Range("I126").Value = 795.37036
$ws.Range("J126").Value = 8336195.5
$ws.Range("K126").Value = 2386.11108
$ws.Range("L126").Value = 25008586.5
$ws.Range("M126").Value = 83.88891999999987
$ws.Range("N126").Value = -25013526.5
$ws.Range("H132").Value = 8366.076999999999
$ws.Range("I132").Value = 1895.46
$ws.Range("J132").Value = 19920.75
$ws.Range("K132").Value = 5686.38
$ws.Range("L132").Value = 59762.25
$ws.Range("M132").Value = -3156.38
$ws.Range("N132").Value = -64822.25
$ws.Range("H136").Value = 895.6177
$ws.Range("I136").Value = 608.24
$ws.Range("J136").Value = 1693.8889
$ws.Range("K136").Value = 1824.72
$ws.Range("L136").Value = 5081.6667
$ws.Range("M136").Value = 725.28
$ws.Range("N136").Value = -10181.6667
